$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A gets a custom width (~17.13 chars); other columns stay default ---
$ws.Columns.Item(1).ColumnWidth = 16.42

# --- Row 2: extend the title cell's formatting across A2:O2 and merge + center it ---
$titleRange = $ws.Range("A2:O2")
$titleRange.Merge()
$titleRange.HorizontalAlignment = -4108   # xlCenter
$titleRange.VerticalAlignment = -4108     # xlCenter

# --- Rows 7-17: bump the row height slightly (21.6 -> 21.75), matching row 6 ---
for ($r = 7; $r -le 17; $r++) {
    $ws.Rows.Item($r).RowHeight = 21.75
}

# --- Header / footer: fix the font-face label (标准 -> Regular) ---
$ws.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ws.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12页 &P'

# --- Move the active selection to M16 ---
$ws.Range("M16").Select() | Out-Null
